# Included support for zalenium
# Set the "Execute" flag for the SigninPageTests error-message test case to "No"
# and update the active selection on the RunManager sheet to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunManager")

$ws.Range("C4").Value = "No"

$ws.Range("C4").Select()
